# "update scripts wuth new tpm" - recompute the NATMI ligand/receptor
# edge-weight table (Wnt9a-Fzd4) for the new TPM-based expression values.
# Columns G-T of rows 2-10 are refreshed; A-F are untouched identifiers/keys.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  "G2" = 0.2304126666666667
  "H2" = 0.691238
  "I2" = 0.03265479005310033
  "J2" = 0.03265479005310033
  "K2" = 3
  "L2" = 1
  "M2" = 26.532132
  "N2" = 79.596396
  "O2" = 0.3960736634233649
  "P2" = 0.3960736634233648
  "Q2" = 6.113339286472001
  "R2" = 55.020053578248
  "S2" = 0.0129337023246523
  "T2" = 0.0129337023246523
  "G3" = 0.2304126666666667
  "H3" = 0.691238
  "I3" = 0.03265479005310033
  "J3" = 0.03265479005310033
  "O3" = 0.2505213219764053
  "P3" = 0.2505213219764053
  "Q3" = 3.866760103410889
  "R3" = 34.800840930698
  "S3" = 0.008180721172964664
  "T3" = 0.008180721172964664
  "G4" = 0.2304126666666667
  "H4" = 0.691238
  "I4" = 0.03265479005310033
  "J4" = 0.03265479005310033
  "M4" = 23.67385
  "N4" = 71.02154999999999
  "O4" = 0.3534050146002298
  "P4" = 0.3534050146002298
  "Q4" = 5.454754908766667
  "R4" = 49.09279417889999
  "S4" = 0.01154036655548336
  "T4" = 0.01154036655548336
  "I5" = 0.8526166070240881
  "J5" = 0.8526166070240883
  "K5" = 3
  "L5" = 1
  "M5" = 26.532132
  "N5" = 79.596396
  "O5" = 0.3960736634233649
  "P5" = 0.3960736634233648
  "Q5" = 159.6192960219
  "R5" = 1436.5736641971
  "S5" = 0.33769898303963
  "T5" = 0.33769898303963
  "I6" = 0.8526166070240881
  "J6" = 0.8526166070240883
  "O6" = 0.2505213219764053
  "P6" = 0.2505213219764053
  "R6" = 908.649997983975
  "S6" = 0.2135986395307118
  "T6" = 0.2135986395307118
  "I7" = 0.8526166070240881
  "J7" = 0.8526166070240883
  "M7" = 23.67385
  "N7" = 71.02154999999999
  "O7" = 0.3534050146002298
  "P7" = 0.3534050146002298
  "Q7" = 142.42365713875
  "R7" = 1281.81291424875
  "S7" = 0.3013189844537463
  "T7" = 0.3013189844537463
  "G8" = 0.8095266666666667
  "H8" = 2.42858
  "I8" = 0.1147286029228115
  "J8" = 0.1147286029228115
  "K8" = 3
  "L8" = 1
  "M8" = 26.532132
  "N8" = 79.596396
  "O8" = 0.3960736634233649
  "P8" = 0.3960736634233648
  "Q8" = 21.47846837752
  "R8" = 193.30621539768
  "S8" = 0.04544097805908252
  "T8" = 0.04544097805908252
  "G9" = 0.8095266666666667
  "H9" = 2.42858
  "I9" = 0.1147286029228115
  "J9" = 0.1147286029228115
  "O9" = 0.2505213219764053
  "P9" = 0.2505213219764053
  "Q9" = 13.58538774190889
  "R9" = 122.26848967718
  "S9" = 0.02874196127272882
  "T9" = 0.02874196127272882
  "G10" = 0.8095266666666667
  "H10" = 2.42858
  "I10" = 0.1147286029228115
  "J10" = 0.1147286029228115
  "M10" = 23.67385
  "N10" = 71.02154999999999
  "O10" = 0.3534050146002298
  "P10" = 0.3534050146002298
  "Q10" = 19.16461287766667  # = G10 * M10 (edge avg weight), derived consistently with R10/S10/T10
  "R10" = 172.481515899
  "S10" = 0.04054566359100017
  "T10" = 0.04054566359100018
}

foreach ($cellRef in $updates.Keys) {
  $ws.Range($cellRef).Value = $updates[$cellRef]
}
